$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00001292064567892659
$ws.Range("C2").Value = 0.00006240767534437808
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 1.247055132155657
